$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-10) currently holds the literal year 2023 as a number.
# Replace it with the text "31/12/2023" and give it a short-date number
# format (built-in numFmtId 14).
$ws.Cells.Item(2, 3).Value = "31/12/2023"
$ws.Cells.Item(2, 3).NumberFormat = "mm-dd-yy"

# Copy the newly formatted cell's style onto the rest of the column so every
# row shares the same cell-format record instead of minting a new one each
# time (mirrors how Excel reuses an existing xf when you fill/copy a format).
$ws.Cells.Item(2, 3).Copy()
$fillRange = $ws.Range("C3:C10")
$fillRange.PasteSpecial(-4122)  # xlPasteFormats
$fillRange.Value = "31/12/2023"
$excel.CutCopyMode = $false

# Reset the sheet selection back to the top-left cell.
$ws.Range("A1").Select()
